# Q4 2022 Fiscal Data update
# Appends the new RMO No. 53-2022 ("Revised CY2022 RMO Goal Annexes",
# issued Dec 7, 2022) row to the BIR Collection Goals table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new text cells first, in the same left-to-right order the
# shared-string table picks up new, not-yet-seen strings (B, C, D, E, A)
# so new <si> entries land at the same indices as the authored workbook.
$ws.Range("B37").Value = "RMO No. 53-2022"
$ws.Range("C37").Value = "December 7, 2022"
$ws.Range("D37").Value = "2022"
$ws.Range("E37").Value = "Final"
$ws.Range("A37").Value = "https://www.bir.gov.ph/images/bir_files/internal_communications_3/2022/Full%20Text/Revised%20CY2022%20RMO%20Goal%20Annexes.pdf"

# Numeric goal columns for the new row.
$ws.Range("F37").Value = 2392587
$ws.Range("G37").Value = 1197966
$ws.Range("H37").Value = 365197
$ws.Range("I37").Value = 430160.553
$ws.Range("J37").Value = 153695.364
$ws.Range("K37").Value = 245568

# Move the frozen-pane selection to just below the new last row, matching
# the saved view state of the updated workbook.
$ws.Range("A38").Select()
